$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A63").Value = "11/03/2025"
$ws.Range("B63").Value = 0.2058740788361836
$ws.Range("C63").Value = 0.7941259211638164
